# Update "horarios 141" schedule workbook with the latest scrape (04:52:00).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 04:52:00"

$data1 = @(
    @("04:52:00", "04:53", "11_ETCHEVERRY",     1,   "LP1912"),
    @("04:52:00", "05:12", "17_ROMERO",         20,  "LP1912"),
    @("04:52:00", "05:22", "23_HERNANDEZ",      30,  "LP1912"),
    @("04:52:00", "05:32", "81_EL PELIGRO",     40,  "LP1912"),
    @("04:52:00", "05:44", "14_ABASTO",         52,  "LP1912"),
    @("04:52:00", "05:52", "17_ROMERO",         60,  "LP1912"),
    @("04:52:00", "06:01", "16_SANTA ANA",      69,  "LP1912"),
    @("04:52:00", "06:04", "10_OLMOS",          72,  "LP1912"),
    @("04:52:00", "06:11", "215A_EL PATO",      79,  "LP1912"),
    @("04:52:00", "06:24", "11_ETCHEVERRY",     92,  "LP1912"),
    @("04:52:00", "06:27", "23_HERNANDEZ",      95,  "LP1912"),
    @("04:52:00", "06:31", "17X38_ROMERO",      99,  "LP1912"),
    @("04:52:00", "06:31", "16_SANTA ANA",      99,  "LP1912"),
    @("04:52:00", "06:39", "225_C ROCA-H SUR",  107, "LP1912"),
    @("04:52:00", "06:51", "215A_EL PATO",      119, "LP1912")
)

$row = 6
foreach ($r in $data1) {
    $ws1.Cells.Item($row, 1).Value = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Cells.Item($row, 4).Value = $r[3]
    $ws1.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 04:52:00"

$data2 = @(
    @("04:52:00", "06:11", "215A_EL PATO", 79,  "LP1912"),
    @("04:52:00", "06:51", "215A_EL PATO", 119, "LP1912")
)

$row = 6
foreach ($r in $data2) {
    $ws2.Cells.Item($row, 1).Value = $r[0]
    $ws2.Cells.Item($row, 2).Value = $r[1]
    $ws2.Cells.Item($row, 3).Value = $r[2]
    $ws2.Cells.Item($row, 4).Value = $r[3]
    $ws2.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 04:52:00"
